$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26, shifting the existing rows 26-35 down to 27-36.
$ws.Rows("26:26").Insert()

# Populate the new row 26 with the same constant columns used throughout this
# data block, plus the specific values for this new record.
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 44523
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 300000000
$ws.Range("G26").Value = "Espárragos"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 520
$ws.Range("K26").Value = 1800
$ws.Range("L26").Value = 1800
$ws.Range("M26").Value = 1800
$ws.Range("N26").Value = "$/kilo"
$ws.Range("O26").Value = "Provincia de Linares"
$ws.Range("P26").Value = 1800
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"

# Match the date cell style used by the rest of the column (style index 2).
$ws.Range("D26").NumberFormat = $ws.Range("D27").NumberFormat
